# Stock updated by raj time 1:17
# Apply updated stock values (column B "qty-ish" figures and their
# corresponding column D totals) on the "PATRIKA 24-25" sheet, plus a
# small width tweak on column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PATRIKA 24-25")

# Narrow column D slightly (9.42578125 -> 8.42578125 char-width units).
# Excel's ColumnWidth COM setter quantizes to whole-pixel steps (MDW-7
# grid), so we pick the character width that lands on the closest
# reachable pixel column to the target.
$ws.Columns.Item(4).ColumnWidth = 7.6

# Row -> (B value, D value)
$updates = @(
    @{ Row = 9;   B = 722.5;  D = 722.5 }
    @{ Row = 10;  B = 304;    D = 304 }
    @{ Row = 44;  B = 150.5;  D = 338.63 }
    @{ Row = 50;  B = 133.5;  D = 146.85 }
    @{ Row = 52;  B = 161;    D = 177.1 }
    @{ Row = 53;  B = 48;     D = 52.8 }
    @{ Row = 54;  B = 39;     D = 42.9 }
    @{ Row = 78;  B = 72.5 }
    @{ Row = 154; B = 98;     D = 289.10000000000002 }
    @{ Row = 155; B = 80;     D = 236 }
    @{ Row = 202; B = 85;     D = 349.18 }
    @{ Row = 212; B = 36;     D = 189 }
    @{ Row = 227; B = 10;     D = 57 }
    @{ Row = 233; B = 7.5;    D = 30 }
    @{ Row = 241; B = 189;    D = 803.25 }
    @{ Row = 308; B = 27.7;   D = 243.76 }
    @{ Row = 417; B = 21.95;  D = 69.14 }
    @{ Row = 502; B = 44;     D = 140.80000000000001 }
    @{ Row = 511; B = 59.5;   D = 198.14 }
    @{ Row = 568; B = 70;     D = 210 }
    @{ Row = 580; B = 125;    D = 348.92 }
    @{ Row = 585; B = 147;    D = 199.5 }
    @{ Row = 596; B = 217;    D = 112.84 }
    @{ Row = 597; B = 410;    D = 213.2 }
    @{ Row = 604; B = 34029.47; D = 99898.93 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 2).Value = $u.B
    if ($u.ContainsKey('D')) {
        $ws.Cells.Item($r, 4).Value = $u.D
    }
}
